$wb = $excel.ActiveWorkbook

# Remove the leftover custom-format styling on S4/S5/S6 rows 20 and 25
foreach ($name in @("S4","S5","S6")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Rows("20:20").ClearFormats()
    $ws.Rows("25:25").ClearFormats()
}

# Add the new S4_SU (substance use) sheet after S6
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "S4_SU"

$ws.Cells.Item(1,1).Value = 'question'
$ws.Cells.Item(1,2).Value = 'routing'

$ws.Cells.Item(78,2).Value = 'drugs_use_frequency != "I take drugs at least once a week or more" & drugs_use_frequency != "I take drugs once or twice a month"'
$ws.Cells.Item(78,1).Value = 'drugs_use_last_year_type_other'
$ws.Cells.Item(77,2).Value = 'drugs_use_frequency != "I take drugs at least once a week or more" & drugs_use_frequency != "I take drugs once or twice a month"'
$ws.Cells.Item(77,1).Value = 'drugs_use_last_year_type_tanning_pills'
$ws.Cells.Item(76,2).Value = 'drugs_use_frequency != "I take drugs at least once a week or more" & drugs_use_frequency != "I take drugs once or twice a month"'
$ws.Cells.Item(76,1).Value = 'drugs_use_last_year_type_diet_pills'
$ws.Cells.Item(75,2).Value = 'drugs_use_frequency != "I take drugs at least once a week or more" & drugs_use_frequency != "I take drugs once or twice a month"'
$ws.Cells.Item(75,1).Value = 'drugs_use_last_year_type_2c'
$ws.Cells.Item(74,2).Value = 'drugs_use_frequency != "I take drugs at least once a week or more" & drugs_use_frequency != "I take drugs once or twice a month"'
$ws.Cells.Item(74,1).Value = 'drugs_use_last_year_type_lsd'
$ws.Cells.Item(73,2).Value = 'drugs_use_frequency != "I take drugs at least once a week or more" & drugs_use_frequency != "I take drugs once or twice a month"'
$ws.Cells.Item(73,1).Value = 'drugs_use_last_year_type_synthetic_cannabinoids'
$ws.Cells.Item(72,2).Value = 'drugs_use_frequency != "I take drugs at least once a week or more" & drugs_use_frequency != "I take drugs once or twice a month"'
$ws.Cells.Item(72,1).Value = 'drugs_use_last_year_type_ketamine'
$ws.Cells.Item(71,2).Value = 'drugs_use_frequency != "I take drugs at least once a week or more" & drugs_use_frequency != "I take drugs once or twice a month"'
$ws.Cells.Item(71,1).Value = 'drugs_use_last_year_type_gear'
$ws.Cells.Item(70,2).Value = 'drugs_use_frequency != "I take drugs at least once a week or more" & drugs_use_frequency != "I take drugs once or twice a month"'
$ws.Cells.Item(70,1).Value = 'drugs_use_last_year_type_cannabis'
$ws.Cells.Item(69,2).Value = 'drugs_use_frequency != "I take drugs at least once a week or more" & drugs_use_frequency != "I take drugs once or twice a month"'
$ws.Cells.Item(69,1).Value = 'drugs_use_last_year_type_steroids'
$ws.Cells.Item(68,2).Value = 'drugs_use_frequency != "I take drugs at least once a week or more" & drugs_use_frequency != "I take drugs once or twice a month"'
$ws.Cells.Item(68,1).Value = 'drugs_use_last_year_type_cocaine'
$ws.Cells.Item(67,2).Value = 'drugs_use_frequency != "I take drugs at least once a week or more" & drugs_use_frequency != "I take drugs once or twice a month"'
$ws.Cells.Item(67,1).Value = 'drugs_use_last_year_type_mdma'
$ws.Cells.Item(66,2).Value = 'drugs_use_frequency != "I take drugs at least once a week or more" & drugs_use_frequency != "I take drugs once or twice a month"'
$ws.Cells.Item(66,1).Value = 'drugs_use_last_year_type_methadone'
$ws.Cells.Item(65,2).Value = 'drugs_use_frequency != "I take drugs at least once a week or more" & drugs_use_frequency != "I take drugs once or twice a month"'
$ws.Cells.Item(65,1).Value = 'drugs_use_last_year_type_mushrooms'
$ws.Cells.Item(64,2).Value = 'drugs_use_frequency != "I take drugs at least once a week or more" & drugs_use_frequency != "I take drugs once or twice a month"'
$ws.Cells.Item(64,1).Value = 'drugs_use_last_year_type_heroin'
$ws.Cells.Item(63,2).Value = 'drugs_use_frequency != "I take drugs at least once a week or more" & drugs_use_frequency != "I take drugs once or twice a month"'
$ws.Cells.Item(63,1).Value = 'drugs_use_last_year_type_benzos'
$ws.Cells.Item(62,2).Value = 'drugs_use_frequency != "I take drugs at least once a week or more" & drugs_use_frequency != "I take drugs once or twice a month"'
$ws.Cells.Item(62,1).Value = 'drugs_use_last_year_type_cyroban'
$ws.Cells.Item(61,2).Value = 'drugs_use_frequency != "I take drugs at least once a week or more" & drugs_use_frequency != "I take drugs once or twice a month"'
$ws.Cells.Item(61,1).Value = 'drugs_use_last_year_type_ecstasy'
$ws.Cells.Item(60,2).Value = 'drugs_use_frequency != "I take drugs at least once a week or more" & drugs_use_frequency != "I take drugs once or twice a month"'
$ws.Cells.Item(60,1).Value = 'drugs_use_last_year_type_amphetamines'
$ws.Cells.Item(59,2).Value = 'drugs_use_frequency != "I take drugs at least once a week or more" & drugs_use_frequency != "I take drugs once or twice a month"'
$ws.Cells.Item(59,1).Value = 'drugs_use_last_year_type_solvents'
$ws.Cells.Item(58,2).Value = 'drugs_use_frequency != "I take drugs at least once a week or more" & drugs_use_frequency != "I take drugs once or twice a month"'
$ws.Cells.Item(58,1).Value = 'drugs_use_last_year_type_cannabis'
$ws.Cells.Item(57,2).Value = 'drugs_ever_taken != "Yes"'
$ws.Cells.Item(57,1).Value = 'drugs_use_last_year'
$ws.Cells.Item(56,2).Value = 'drugs_ever_taken != "Yes"'
$ws.Cells.Item(56,1).Value = 'drugs_use_frequency'
$ws.Cells.Item(55,2).Value = 'alcohol_ever_had_any != "Yes"'
$ws.Cells.Item(55,1).Value = 'alcohol_last_provided_by'
$ws.Cells.Item(54,2).Value = 'alcohol_ever_had_any != "Yes"'
$ws.Cells.Item(54,1).Value = 'alcohol_drinking_allowed_at_home'
$ws.Cells.Item(53,2).Value = 'alcohol_ever_had_any != "Yes"'
$ws.Cells.Item(53,1).Value = 'alcohol_usual_drinking_location_elsewhere'
$ws.Cells.Item(52,2).Value = 'alcohol_ever_had_any != "Yes"'
$ws.Cells.Item(52,1).Value = 'alcohol_usual_drinking_location_outdoors'
$ws.Cells.Item(51,2).Value = 'alcohol_ever_had_any != "Yes"'
$ws.Cells.Item(51,1).Value = 'alcohol_usual_drinking_location_someones_home'
$ws.Cells.Item(50,2).Value = 'alcohol_ever_had_any != "Yes"'
$ws.Cells.Item(50,1).Value = 'alcohol_usual_drinking_location_home'
$ws.Cells.Item(49,2).Value = 'alcohol_ever_had_any != "Yes"'
$ws.Cells.Item(49,1).Value = 'alcohol_usual_drinking_location_party'
$ws.Cells.Item(48,2).Value = 'alcohol_ever_had_any != "Yes"'
$ws.Cells.Item(48,1).Value = 'alcohol_usual_drinking_location_club'
$ws.Cells.Item(47,2).Value = 'alcohol_ever_had_any != "Yes"'
$ws.Cells.Item(47,1).Value = 'alcohol_usual_drinking_location_pub'
$ws.Cells.Item(46,2).Value = 'alcohol_ever_had_any != "Yes"'
$ws.Cells.Item(46,1).Value = 'alcohol_usual_source'
$ws.Cells.Item(45,2).Value = 'alcohol_ever_had_any != "Yes"'
$ws.Cells.Item(45,1).Value = 'alcohol_frequency_getting_drunk'
$ws.Cells.Item(44,2).Value = 'alcohol_ever_had_any != "Yes"'
$ws.Cells.Item(44,1).Value = 'alcohol_frequency_type_other'
$ws.Cells.Item(43,2).Value = 'alcohol_ever_had_any != "Yes"'
$ws.Cells.Item(43,1).Value = 'alcohol_frequency_type_wine_fortified'
$ws.Cells.Item(42,2).Value = 'alcohol_ever_had_any != "Yes"'
$ws.Cells.Item(42,1).Value = 'alcohol_frequency_type_cider'
$ws.Cells.Item(41,2).Value = 'alcohol_ever_had_any != "Yes"'
$ws.Cells.Item(41,1).Value = 'alcohol_frequency_type_spirits'
$ws.Cells.Item(40,2).Value = 'alcohol_ever_had_any != "Yes"'
$ws.Cells.Item(40,1).Value = 'alcohol_frequency_type_alcopops'
$ws.Cells.Item(39,2).Value = 'alcohol_ever_had_any != "Yes"'
$ws.Cells.Item(39,1).Value = 'alcohol_frequency_type_wine'
$ws.Cells.Item(38,2).Value = 'alcohol_ever_had_any != "Yes"'
$ws.Cells.Item(38,1).Value = 'alcohol_frequency_type_beer'
$ws.Cells.Item(37,2).Value = 'e_cigarettes_use_frequency != "I use e-cigarettes / vapes once a week or more" & e_cigarettes_use_frequency != "I use e-cigarettes / vapes sometimes, but no more than once a month"'
$ws.Cells.Item(37,1).Value = 'e_cigarettes_source_other'
$ws.Cells.Item(36,2).Value = 'e_cigarettes_use_frequency != "I use e-cigarettes / vapes once a week or more" & e_cigarettes_use_frequency != "I use e-cigarettes / vapes sometimes, but no more than once a month"'
$ws.Cells.Item(36,1).Value = 'e_cigarettes_source_take_without_asking'
$ws.Cells.Item(35,2).Value = 'e_cigarettes_use_frequency != "I use e-cigarettes / vapes once a week or more" & e_cigarettes_use_frequency != "I use e-cigarettes / vapes sometimes, but no more than once a month"'
$ws.Cells.Item(35,1).Value = 'e_cigarettes_source_parents_provide'
$ws.Cells.Item(34,2).Value = 'e_cigarettes_use_frequency != "I use e-cigarettes / vapes once a week or more" & e_cigarettes_use_frequency != "I use e-cigarettes / vapes sometimes, but no more than once a month"'
$ws.Cells.Item(34,1).Value = 'e_cigarettes_source_siblings_provide'
$ws.Cells.Item(33,2).Value = 'e_cigarettes_use_frequency != "I use e-cigarettes / vapes once a week or more" & e_cigarettes_use_frequency != "I use e-cigarettes / vapes sometimes, but no more than once a month"'
$ws.Cells.Item(33,1).Value = 'e_cigarettes_source_friends_provide'
$ws.Cells.Item(32,2).Value = 'e_cigarettes_use_frequency != "I use e-cigarettes / vapes once a week or more" & e_cigarettes_use_frequency != "I use e-cigarettes / vapes sometimes, but no more than once a month"'
$ws.Cells.Item(32,1).Value = 'e_cigarettes_source_ask_adult_unknown'
$ws.Cells.Item(31,2).Value = 'e_cigarettes_use_frequency != "I use e-cigarettes / vapes once a week or more" & e_cigarettes_use_frequency != "I use e-cigarettes / vapes sometimes, but no more than once a month"'
$ws.Cells.Item(31,1).Value = 'e_cigarettes_source_ask_adult_known'
$ws.Cells.Item(30,2).Value = 'e_cigarettes_use_frequency != "I use e-cigarettes / vapes once a week or more" & e_cigarettes_use_frequency != "I use e-cigarettes / vapes sometimes, but no more than once a month"'
$ws.Cells.Item(30,1).Value = 'e_cigarettes_source_ask_minor_known'
$ws.Cells.Item(29,2).Value = 'e_cigarettes_use_frequency != "I use e-cigarettes / vapes once a week or more" & e_cigarettes_use_frequency != "I use e-cigarettes / vapes sometimes, but no more than once a month"'
$ws.Cells.Item(29,1).Value = 'e_cigarettes_source_someone_else'
$ws.Cells.Item(28,2).Value = 'e_cigarettes_use_frequency != "I use e-cigarettes / vapes once a week or more" & e_cigarettes_use_frequency != "I use e-cigarettes / vapes sometimes, but no more than once a month"'
$ws.Cells.Item(28,1).Value = 'e_cigarettes_source_friends_relatives'
$ws.Cells.Item(27,2).Value = 'e_cigarettes_use_frequency != "I use e-cigarettes / vapes once a week or more" & e_cigarettes_use_frequency != "I use e-cigarettes / vapes sometimes, but no more than once a month"'
$ws.Cells.Item(27,1).Value = 'e_cigarettes_source_internet'
$ws.Cells.Item(26,2).Value = 'e_cigarettes_use_frequency != "I use e-cigarettes / vapes once a week or more" & e_cigarettes_use_frequency != "I use e-cigarettes / vapes sometimes, but no more than once a month"'
$ws.Cells.Item(26,1).Value = 'e_cigarettes_source_street_market'
$ws.Cells.Item(25,2).Value = 'e_cigarettes_use_frequency != "I use e-cigarettes / vapes once a week or more" & e_cigarettes_use_frequency != "I use e-cigarettes / vapes sometimes, but no more than once a month"'
$ws.Cells.Item(25,1).Value = 'e_cigarettes_source_hop_other'
$ws.Cells.Item(24,2).Value = 'e_cigarettes_use_frequency != "I use e-cigarettes / vapes once a week or more" & e_cigarettes_use_frequency != "I use e-cigarettes / vapes sometimes, but no more than once a month"'
$ws.Cells.Item(24,1).Value = 'e_cigarettes_source_van'
$ws.Cells.Item(23,2).Value = 'e_cigarettes_use_frequency != "I use e-cigarettes / vapes once a week or more" & e_cigarettes_use_frequency != "I use e-cigarettes / vapes sometimes, but no more than once a month"'
$ws.Cells.Item(23,1).Value = 'e_cigarettes_source_garage_shop'
$ws.Cells.Item(22,2).Value = 'e_cigarettes_use_frequency != "I use e-cigarettes / vapes once a week or more" & e_cigarettes_use_frequency != "I use e-cigarettes / vapes sometimes, but no more than once a month"'
$ws.Cells.Item(22,1).Value = 'e_cigarettes_source_newsagent'
$ws.Cells.Item(21,2).Value = 'e_cigarettes_use_frequency != "I use e-cigarettes / vapes once a week or more" & e_cigarettes_use_frequency != "I use e-cigarettes / vapes sometimes, but no more than once a month"'
$ws.Cells.Item(21,1).Value = 'e_cigarettes_source_supermarket'
$ws.Cells.Item(20,2).Value = 'cigarettes_smoking_status != "I usually smoke more than six cigarettes a week" & cigarettes_smoking_status != "I usually smoke between one and six cigarettes a week" & cigarettes_smoking_status != "I sometimes smoke cigarettes now but I don''t smoke as many as one a week"'
$ws.Cells.Item(20,1).Value = 'cigarettes_daily_number'
$ws.Cells.Item(19,2).Value = 'cigarettes_smoking_status != "I usually smoke more than six cigarettes a week" & cigarettes_smoking_status != "I usually smoke between one and six cigarettes a week" & cigarettes_smoking_status != "I sometimes smoke cigarettes now but I don''t smoke as many as one a week"'
$ws.Cells.Item(19,1).Value = 'cigarettes_attempts_buying'
$ws.Cells.Item(18,2).Value = 'cigarettes_smoking_status != "I usually smoke more than six cigarettes a week" & cigarettes_smoking_status != "I usually smoke between one and six cigarettes a week" & cigarettes_smoking_status != "I sometimes smoke cigarettes now but I don''t smoke as many as one a week"'
$ws.Cells.Item(18,1).Value = 'cigarettes_source_other'
$ws.Cells.Item(17,2).Value = 'cigarettes_smoking_status != "I usually smoke more than six cigarettes a week" & cigarettes_smoking_status != "I usually smoke between one and six cigarettes a week" & cigarettes_smoking_status != "I sometimes smoke cigarettes now but I don''t smoke as many as one a week"'
$ws.Cells.Item(17,1).Value = 'cigarettes_source_take_without_asking'
$ws.Cells.Item(16,2).Value = 'cigarettes_smoking_status != "I usually smoke more than six cigarettes a week" & cigarettes_smoking_status != "I usually smoke between one and six cigarettes a week" & cigarettes_smoking_status != "I sometimes smoke cigarettes now but I don''t smoke as many as one a week"'
$ws.Cells.Item(16,1).Value = 'cigarettes_source_parents_provide'
$ws.Cells.Item(15,2).Value = 'cigarettes_smoking_status != "I usually smoke more than six cigarettes a week" & cigarettes_smoking_status != "I usually smoke between one and six cigarettes a week" & cigarettes_smoking_status != "I sometimes smoke cigarettes now but I don''t smoke as many as one a week"'
$ws.Cells.Item(15,1).Value = 'cigarettes_source_siblings_provide'
$ws.Cells.Item(14,2).Value = 'cigarettes_smoking_status != "I usually smoke more than six cigarettes a week" & cigarettes_smoking_status != "I usually smoke between one and six cigarettes a week" & cigarettes_smoking_status != "I sometimes smoke cigarettes now but I don''t smoke as many as one a week"'
$ws.Cells.Item(14,1).Value = 'cigarettes_source_friends_provide'
$ws.Cells.Item(13,2).Value = 'cigarettes_smoking_status != "I usually smoke more than six cigarettes a week" & cigarettes_smoking_status != "I usually smoke between one and six cigarettes a week" & cigarettes_smoking_status != "I sometimes smoke cigarettes now but I don''t smoke as many as one a week"'
$ws.Cells.Item(13,1).Value = 'cigarettes_source_ask_adult_unknown'
$ws.Cells.Item(12,2).Value = 'cigarettes_smoking_status != "I usually smoke more than six cigarettes a week" & cigarettes_smoking_status != "I usually smoke between one and six cigarettes a week" & cigarettes_smoking_status != "I sometimes smoke cigarettes now but I don''t smoke as many as one a week"'
$ws.Cells.Item(12,1).Value = 'cigarettes_source_ask_adult_known'
$ws.Cells.Item(11,2).Value = 'cigarettes_smoking_status != "I usually smoke more than six cigarettes a week" & cigarettes_smoking_status != "I usually smoke between one and six cigarettes a week" & cigarettes_smoking_status != "I sometimes smoke cigarettes now but I don''t smoke as many as one a week"'
$ws.Cells.Item(11,1).Value = 'cigarettes_source_ask_minor_known'
$ws.Cells.Item(10,2).Value = 'cigarettes_smoking_status != "I usually smoke more than six cigarettes a week" & cigarettes_smoking_status != "I usually smoke between one and six cigarettes a week" & cigarettes_smoking_status != "I sometimes smoke cigarettes now but I don''t smoke as many as one a week"'
$ws.Cells.Item(10,1).Value = 'cigarettes_source_someone_else'
$ws.Cells.Item(9,2).Value = 'cigarettes_smoking_status != "I usually smoke more than six cigarettes a week" & cigarettes_smoking_status != "I usually smoke between one and six cigarettes a week" & cigarettes_smoking_status != "I sometimes smoke cigarettes now but I don''t smoke as many as one a week"'
$ws.Cells.Item(9,1).Value = 'cigarettes_source_friends_relatives'
$ws.Cells.Item(8,2).Value = 'cigarettes_smoking_status != "I usually smoke more than six cigarettes a week" & cigarettes_smoking_status != "I usually smoke between one and six cigarettes a week" & cigarettes_smoking_status != "I sometimes smoke cigarettes now but I don''t smoke as many as one a week"'
$ws.Cells.Item(8,1).Value = 'cigarettes_source_internet'
$ws.Cells.Item(7,2).Value = 'cigarettes_smoking_status != "I usually smoke more than six cigarettes a week" & cigarettes_smoking_status != "I usually smoke between one and six cigarettes a week" & cigarettes_smoking_status != "I sometimes smoke cigarettes now but I don''t smoke as many as one a week"'
$ws.Cells.Item(7,1).Value = 'cigarettes_source_street_market'
$ws.Cells.Item(6,2).Value = 'cigarettes_smoking_status != "I usually smoke more than six cigarettes a week" & cigarettes_smoking_status != "I usually smoke between one and six cigarettes a week" & cigarettes_smoking_status != "I sometimes smoke cigarettes now but I don''t smoke as many as one a week"'
$ws.Cells.Item(6,1).Value = 'cigarettes_source_hop_other'
$ws.Cells.Item(5,2).Value = 'cigarettes_smoking_status != "I usually smoke more than six cigarettes a week" & cigarettes_smoking_status != "I usually smoke between one and six cigarettes a week" & cigarettes_smoking_status != "I sometimes smoke cigarettes now but I don''t smoke as many as one a week"'
$ws.Cells.Item(5,1).Value = 'cigarettes_source_van'
$ws.Cells.Item(4,2).Value = 'cigarettes_smoking_status != "I usually smoke more than six cigarettes a week" & cigarettes_smoking_status != "I usually smoke between one and six cigarettes a week" & cigarettes_smoking_status != "I sometimes smoke cigarettes now but I don''t smoke as many as one a week"'
$ws.Cells.Item(4,1).Value = 'cigarettes_source_garage_shop'
$ws.Cells.Item(3,2).Value = 'cigarettes_smoking_status != "I usually smoke more than six cigarettes a week" & cigarettes_smoking_status != "I usually smoke between one and six cigarettes a week" & cigarettes_smoking_status != "I sometimes smoke cigarettes now but I don''t smoke as many as one a week"'
$ws.Cells.Item(3,1).Value = 'cigarettes_source_newsagent'
$ws.Cells.Item(2,2).Value = 'cigarettes_smoking_status != "I usually smoke more than six cigarettes a week" & cigarettes_smoking_status != "I usually smoke between one and six cigarettes a week" & cigarettes_smoking_status != "I sometimes smoke cigarettes now but I don''t smoke as many as one a week"'
$ws.Cells.Item(2,1).Value = 'cigarettes_source_supermarket'
